$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update existing D column values (rows 2-5 "tes", rows 6-9 "goo")
$ws.Range("D2").Value = 87
$ws.Range("D3").Value = 90
$ws.Range("D4").Value = 89
$ws.Range("D5").Value = 85
$ws.Range("D6").Value = 41
$ws.Range("D7").Value = 44
$ws.Range("D8").Value = 43
$ws.Range("D9").Value = 50

# Copy formatting from the last existing data row (row 9) down to the new rows first
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Add new rows 10-13 for "ms"
$ws.Range("A10").Value = "ms"
$ws.Range("B10").Value = 41244
$ws.Range("C10").Value = "us"
$ws.Range("D10").Value = 20

$ws.Range("A11").Value = "ms"
$ws.Range("B11").Value = 41245
$ws.Range("C11").Value = "us"
$ws.Range("D11").Value = 20

$ws.Range("A12").Value = "ms"
$ws.Range("B12").Value = 41246
$ws.Range("C12").Value = "ch"
$ws.Range("D12").Value = 20

$ws.Range("A13").Value = "ms"
$ws.Range("B13").Value = 41247
$ws.Range("C13").Value = "nz"
$ws.Range("D13").Value = 20

# Update selection to match the new active cell (row below new data)
$ws.Rows.Item(14).Select()
